{"js": "// Replace the date and every two-digit multiplication problem's text with\n// the new values from the commit. Each old string is unique in the\n// document, so a simple search + replace per pair is safe and unambiguous.\nconst replacements = [\n  [\"2024-04-04 Thursday\", \"2024-04-05 Friday\"],\n  [\"78\u00d722=\", \"97\u00d782=\"],\n  [\"56\u00d726=\", \"66\u00d798=\"],\n  [\"68\u00d782=\", \"65\u00d774=\"],\n  [\"60\u00d724=\", \"68\u00d718=\"],\n  [\"88\u00d757=\", \"43\u00d775=\"],\n  [\"25\u00d747=\", \"62\u00d718=\"],\n  [\"11\u00d778=\", \"64\u00d714=\"],\n  [\"62\u00d752=\", \"23\u00d760=\"],\n  [\"34\u00d791=\", \"96\u00d759=\"],\n  [\"76\u00d780=\", \"71\u00d761=\"],\n  [\"33\u00d741=\", \"29\u00d726=\"],\n  [\"79\u00d759=\", \"56\u00d777=\"],\n  [\"91\u00d772=\", \"18\u00d790=\"],\n  [\"73\u00d763=\", \"85\u00d743=\"],\n  [\"68\u00d757=\", \"47\u00d726=\"],\n  [\"60\u00d740=\", \"22\u00d735=\"],\n  [\"16\u00d756=\", \"93\u00d727=\"],\n  [\"58\u00d778=\", \"70\u00d713=\"],\n  [\"39\u00d757=\", \"97\u00d771=\"],\n  [\"72\u00d755=\", \"88\u00d765=\"],\n  [\"23\u00d771=\", \"87\u00d783=\"],\n  [\"86\u00d736=\", \"96\u00d711=\"],\n  [\"77\u00d797=\", \"24\u00d716=\"],\n  [\"42\u00d779=\", \"43\u00d743=\"],\n  [\"11\u00d773=\", \"73\u00d749=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and every two-digit multiplication problem's text with\n# the new values from the commit. Each old string is unique in the\n# document, so Find/Replace per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-04-04 Thursday\", \"2024-04-05 Friday\"),\n    @(\"78\u00d722=\", \"97\u00d782=\"),\n    @(\"56\u00d726=\", \"66\u00d798=\"),\n    @(\"68\u00d782=\", \"65\u00d774=\"),\n    @(\"60\u00d724=\", \"68\u00d718=\"),\n    @(\"88\u00d757=\", \"43\u00d775=\"),\n    @(\"25\u00d747=\", \"62\u00d718=\"),\n    @(\"11\u00d778=\", \"64\u00d714=\"),\n    @(\"62\u00d752=\", \"23\u00d760=\"),\n    @(\"34\u00d791=\", \"96\u00d759=\"),\n    @(\"76\u00d780=\", \"71\u00d761=\"),\n    @(\"33\u00d741=\", \"29\u00d726=\"),\n    @(\"79\u00d759=\", \"56\u00d777=\"),\n    @(\"91\u00d772=\", \"18\u00d790=\"),\n    @(\"73\u00d763=\", \"85\u00d743=\"),\n    @(\"68\u00d757=\", \"47\u00d726=\"),\n    @(\"60\u00d740=\", \"22\u00d735=\"),\n    @(\"16\u00d756=\", \"93\u00d727=\"),\n    @(\"58\u00d778=\", \"70\u00d713=\"),\n    @(\"39\u00d757=\", \"97\u00d771=\"),\n    @(\"72\u00d755=\", \"88\u00d765=\"),\n    @(\"23\u00d771=\", \"87\u00d783=\"),\n    @(\"86\u00d736=\", \"96\u00d711=\"),\n    @(\"77\u00d797=\", \"24\u00d716=\"),\n    @(\"42\u00d779=\", \"43\u00d743=\"),\n    @(\"11\u00d773=\", \"73\u00d749=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $find\n    $rng.Find.Replacement.Text = $replace\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 1  # wdFindContinue\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.MatchWildcards = $false\n    $rng.Find.Execute([ref]$find, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$replace, [ref]2) | Out-Null\n}\n"}
